# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" between "2021-Q4" and "总计", with the
#    same column layout/styling as "2021-Q4", holding two fund rows.
# 2. Insert a new row at the top of the "总计" data (row 2) summarising the
#    "2022-Q1" sheet, pushing the existing "2021-Q4" summary row down to
#    row 3.
#
# Note on sheet identity: the existing "总计" worksheet object tracks its
# *tab position*, not a stable identity -- once a new sheet is spliced in
# in front of it, the old handle silently starts resolving to that new
# sheet instead. "总计" is therefore rebuilt (delete + re-add in its
# original slot, repopulated from values captured up front) rather than
# reused in place, which also keeps the three tabs' underlying sheet
# order/identity consistent with a natural "insert a sheet" edit.

$wb = $excel.ActiveWorkbook

$ws2021 = $wb.Worksheets.Item("2021-Q4")
$wsTotalOld = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------------
# Helper: write a value into a cell as literal TEXT (no numeric coercion,
# no leftover number-format styling) -- mirrors how Excel stores numeric-
# looking strings ("001914", "0.13", ...) as text in this workbook.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($Range, $Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.ClearFormats()
}

# Capture "总计"'s existing header + single data row before the sheet is
# torn down.
$totalHeaderB = $wsTotalOld.Range("B1").Value2
$totalHeaderC = $wsTotalOld.Range("C1").Value2
$totalHeaderD = $wsTotalOld.Range("D1").Value2
$total2021B = $wsTotalOld.Range("B2").Value2
$total2021C = $wsTotalOld.Range("C2").Value2
$total2021D = $wsTotalOld.Range("D2").Value2

$wsTotalOld.Delete() | Out-Null

# ---------------------------------------------------------------------------
# 1) Create the "2022-Q1" sheet right after "2021-Q4", then rebuild "总计"
#    right after that -- restoring the original tab order
#    2021-Q4 / 2022-Q1 / 总计.
# ---------------------------------------------------------------------------
$wsNew = $wb.Worksheets.Add($null, $ws2021)
$wsNew.Name = "2022-Q1"

$wsTotal = $wb.Worksheets.Add($null, $wsNew)
$wsTotal.Name = "总计"

# Header row (B1:H1) -- copy formatting from the 2021-Q4 header so the new
# sheet matches its bold/centered/bordered style exactly.
$ws2021.Range("B1:H1").Copy()
$wsNew.Range("B1:H1").PasteSpecial(-4122)

$wsNew.Range("B1").Value = "基金代码"
$wsNew.Range("C1").Value = "基金名称"
$wsNew.Range("D1").Value = "基金规模"
$wsNew.Range("E1").Value = "股票总仓位"
$wsNew.Range("F1").Value = "仓位占比"
$wsNew.Range("G1").Value = "持有市值(亿元)"
$wsNew.Range("H1").Value = "仓位排名"

# Row 2 -- fund #1. Column A carries the same bold/bordered style as the
# source sheet's index column.
$ws2021.Range("A2").Copy()
$wsNew.Range("A2").PasteSpecial(-4122)
$wsNew.Range("A2").Value = 0

Set-TextValue $wsNew.Range("B2") "001914"
Set-TextValue $wsNew.Range("C2") "中信建投聚利混合A"
Set-TextValue $wsNew.Range("D2") "0.13"
Set-TextValue $wsNew.Range("E2") "39.07"
Set-TextValue $wsNew.Range("F2") "2.04"
Set-TextValue $wsNew.Range("G2") "0.0027"
$wsNew.Range("H2").Value = 8

# Row 3 -- fund #2.
$ws2021.Range("A2").Copy()
$wsNew.Range("A3").PasteSpecial(-4122)
$wsNew.Range("A3").Value = 1

Set-TextValue $wsNew.Range("B3") "000041"
Set-TextValue $wsNew.Range("C3") "华夏全球精选股票(QDII)"
Set-TextValue $wsNew.Range("D3") "0.02"
Set-TextValue $wsNew.Range("E3") "39.07"
Set-TextValue $wsNew.Range("F3") "2.04"
Set-TextValue $wsNew.Range("G3") "0.0004"
$wsNew.Range("H3").Value = 8

# ---------------------------------------------------------------------------
# 2) Rebuild "总计": header row, a new row 2 summarising "2022-Q1", and the
#    original "2021-Q4" summary row pushed down to row 3.
# ---------------------------------------------------------------------------
$ws2021.Range("B1:D1").Copy()
$wsTotal.Range("B1:D1").PasteSpecial(-4122)

$wsTotal.Range("B1").Value = $totalHeaderB
$wsTotal.Range("C1").Value = $totalHeaderC
$wsTotal.Range("D1").Value = $totalHeaderD

$ws2021.Range("A2").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)
$wsTotal.Range("A2").Value = 0

Set-TextValue $wsTotal.Range("B2") "2022-Q1"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0

$ws2021.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)
$wsTotal.Range("A3").Value = 1

Set-TextValue $wsTotal.Range("B3") $total2021B
$wsTotal.Range("C3").Value = $total2021C
$wsTotal.Range("D3").Value = $total2021D

Write-Output "edit complete"
